$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 598 (old rows 598-638 shift down
# to 600-640), matching the way the new weekly records were prepended to this
# block of "Repollo" price rows.
$ws.Rows("598:599").Insert()

# Row 598 - new "Primera" record dated 2022-07-04 (serial 44746)
$ws.Range("A598").Value = 3
$ws.Range("B598").Value = "Femacal de La Calera"
$ws.Range("C598").Value = "Coquimbo"
$ws.Range("D598").Value = 44746
$ws.Range("E598").Value = 5
$ws.Range("F598").Value = 100112006
$ws.Range("G598").Value = "Repollo"
$ws.Range("H598").Value = "Crespo record"
$ws.Range("I598").Value = "Primera"
$ws.Range("J598").Value = 1460
$ws.Range("K598").Value = 1300
$ws.Range("L598").Value = 1400
$ws.Range("M598").Value = 1353
$ws.Range("N598").Value = "$/unidad"
$ws.Range("O598").Value = "Provincia de Quillota"
$ws.Range("P598").Value = 1353
$ws.Range("Q598").Value = 1
$ws.Range("R598").Value = "Hortaliza"

# Row 599 - new "Segunda" record dated 2022-07-04 (serial 44746)
$ws.Range("A599").Value = 3
$ws.Range("B599").Value = "Femacal de La Calera"
$ws.Range("C599").Value = "Coquimbo"
$ws.Range("D599").Value = 44746
$ws.Range("E599").Value = 5
$ws.Range("F599").Value = 100112006
$ws.Range("G599").Value = "Repollo"
$ws.Range("H599").Value = "Crespo record"
$ws.Range("I599").Value = "Segunda"
$ws.Range("J599").Value = 1498
$ws.Range("K599").Value = 1000
$ws.Range("L599").Value = 1200
$ws.Range("M599").Value = 1115
$ws.Range("N599").Value = "$/unidad"
$ws.Range("O599").Value = "Provincia de Quillota"
$ws.Range("P599").Value = 1115
$ws.Range("Q599").Value = 1
$ws.Range("R599").Value = "Hortaliza"
